$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the original (default) style of the data cells, then force
# the Price/Volume columns to Text format so Excel does not silently
# convert numeric-looking strings (e.g. "30.694.63", "1.0000") into
# floating point numbers when we assign them. Restore the original style
# afterwards so no stray style index is left referenced on the cells.
$origStyle = $ws.Range("B2").Style
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "30.694.63"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "1.887.57"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "247.86"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").Value = "0.4741"
$ws.Range("E7").Value = "  -0.36%  "
$ws.Range("D8").Value = "0.2924"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "0.06528"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "21.98"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").Value = "0.07793"
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("B12").Value = "Litecoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D12").Value = "97.02"
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.891.94"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("D14").Value = "0.7361"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").Value = "5.249"
$ws.Range("E15").Value = "  +1.97%  "
$ws.Range("D16").Value = "283.52"
$ws.Range("E16").Value = "  +3.38%  "
$ws.Range("D17").Value = "30.774.52"
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").Value = "13.17"
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("D19").Value = "0.000007549"
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "2.140.98"
$ws.Range("D22").Value = "5.312"
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "6.262"
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").Value = "9.218"
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("D26").Value = "164.24"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").Value = "18.91"
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").Value = "1.920"
$ws.Range("E28").Value = "  -1.31%  "
$ws.Range("D29").Value = "1.341"
$ws.Range("E29").Value = "  -1.87%  "
$ws.Range("D30").Value = "0.09713"
$ws.Range("E30").Value = "  -3.46%  "
$ws.Range("D31").Value = "1.494"
$ws.Range("E31").Value = "  -1.89%  "
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("D33").Value = "4.196"
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("D34").Value = "0.04845"
$ws.Range("E34").Value = "  +0.43%  "
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").Value = "0.6971"
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("D37").Value = "2.723"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").Value = "0.01908"
$ws.Range("E38").Value = "  +2.01%  "
$ws.Range("E39").Value = "  +2.15%  "
$ws.Range("D40").Value = "6.347"
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("D41").Value = "75.96"
$ws.Range("E41").Value = "  +6.42%  "
$ws.Range("D42").Value = "2.014"
$ws.Range("E42").Value = "  +1.02%  "
$ws.Range("D43").Value = "0.4256"
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("D44").Value = "0.9999"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "0.8353"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("D46").Value = "101.29"
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("D47").Value = "9.510"
$ws.Range("E47").Value = "  +1.83%  "
$ws.Range("D48").Value = "7.041"
$ws.Range("E48").Value = "  -0.79%  "
$ws.Range("D49").Value = "35.61"
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").Value = "916.77"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").Value = "0.05752"
$ws.Range("E51").Value = "  +1.92%  "

$dataRange.Style = $origStyle
